$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = [double]"3"
$ws.Cells.Item(2, 6).Value = [double]"1"
$ws.Cells.Item(2, 7).Value = [double]"0.907169"
$ws.Cells.Item(2, 8).Value = [double]"2.721507"
$ws.Cells.Item(2, 9).Value = [double]"0.001369063862079057"
$ws.Cells.Item(2, 10).Value = [double]"0.001369063862079057"
$ws.Cells.Item(2, 13).Value = [double]"0.106124"
$ws.Cells.Item(2, 14).Value = [double]"0.318372"
$ws.Cells.Item(2, 15).Value = [double]"0.08094716512538251"
$ws.Cells.Item(2, 16).Value = [double]"0.08094716512538253"
$ws.Cells.Item(2, 17).Value = [double]"0.09627240295599999"
$ws.Cells.Item(2, 18).Value = [double]"0.866451626604"
$ws.Cells.Item(2, 19).Value = [double]"0.0001108218385109074"
$ws.Cells.Item(2, 20).Value = [double]"0.0001108218385109074"
$ws.Cells.Item(3, 5).Value = [double]"3"
$ws.Cells.Item(3, 6).Value = [double]"1"
$ws.Cells.Item(3, 7).Value = [double]"0.907169"
$ws.Cells.Item(3, 8).Value = [double]"2.721507"
$ws.Cells.Item(3, 9).Value = [double]"0.001369063862079057"
$ws.Cells.Item(3, 10).Value = [double]"0.001369063862079057"
$ws.Cells.Item(3, 15).Value = [double]"0.8331551016962769"
$ws.Cells.Item(3, 16).Value = [double]"0.833155101696277"
$ws.Cells.Item(3, 17).Value = [double]"0.9908913246203332"
$ws.Cells.Item(3, 18).Value = [double]"8.918021921583"
$ws.Cells.Item(3, 19).Value = [double]"0.001140642541239175"
$ws.Cells.Item(3, 20).Value = [double]"0.001140642541239175"
$ws.Cells.Item(4, 5).Value = [double]"3"
$ws.Cells.Item(4, 6).Value = [double]"1"
$ws.Cells.Item(4, 7).Value = [double]"0.907169"
$ws.Cells.Item(4, 8).Value = [double]"2.721507"
$ws.Cells.Item(4, 9).Value = [double]"0.001369063862079057"
$ws.Cells.Item(4, 10).Value = [double]"0.001369063862079057"
$ws.Cells.Item(4, 13).Value = [double]"0.1126143333333333"
$ws.Cells.Item(4, 14).Value = [double]"0.337843"
$ws.Cells.Item(4, 15).Value = [double]"0.08589773317834044"
$ws.Cells.Item(4, 16).Value = [double]"0.08589773317834046"
$ws.Cells.Item(4, 17).Value = [double]"0.1021602321556667"
$ws.Cells.Item(4, 18).Value = [double]"0.919442089401"
$ws.Cells.Item(4, 19).Value = [double]"0.0001175994823289751"
$ws.Cells.Item(4, 20).Value = [double]"0.0001175994823289752"
$ws.Cells.Item(5, 9).Value = [double]"0.002841027838709403"
$ws.Cells.Item(5, 10).Value = [double]"0.002841027838709403"
$ws.Cells.Item(5, 13).Value = [double]"0.106124"
$ws.Cells.Item(5, 14).Value = [double]"0.318372"
$ws.Cells.Item(5, 15).Value = [double]"0.08094716512538251"
$ws.Cells.Item(5, 16).Value = [double]"0.08094716512538253"
$ws.Cells.Item(5, 17).Value = [double]"0.1997807293533333"
$ws.Cells.Item(5, 18).Value = [double]"1.79802656418"
$ws.Cells.Item(5, 19).Value = [double]"0.0002299731495858186"
$ws.Cells.Item(5, 20).Value = [double]"0.0002299731495858187"
$ws.Cells.Item(6, 9).Value = [double]"0.002841027838709403"
$ws.Cells.Item(6, 10).Value = [double]"0.002841027838709403"
$ws.Cells.Item(6, 15).Value = [double]"0.8331551016962769"
$ws.Cells.Item(6, 16).Value = [double]"0.833155101696277"
$ws.Cells.Item(6, 19).Value = [double]"0.002367016837881886"
$ws.Cells.Item(6, 20).Value = [double]"0.002367016837881887"
$ws.Cells.Item(7, 9).Value = [double]"0.002841027838709403"
$ws.Cells.Item(7, 10).Value = [double]"0.002841027838709403"
$ws.Cells.Item(7, 13).Value = [double]"0.1126143333333333"
$ws.Cells.Item(7, 14).Value = [double]"0.337843"
$ws.Cells.Item(7, 15).Value = [double]"0.08589773317834044"
$ws.Cells.Item(7, 16).Value = [double]"0.08589773317834046"
$ws.Cells.Item(7, 17).Value = [double]"0.2119989224772222"
$ws.Cells.Item(7, 18).Value = [double]"1.907990302295"
$ws.Cells.Item(7, 19).Value = [double]"0.0002440378512416975"
$ws.Cells.Item(7, 20).Value = [double]"0.0002440378512416976"
$ws.Cells.Item(8, 7).Value = [double]"84.26343166666668"
$ws.Cells.Item(8, 8).Value = [double]"252.790295"
$ws.Cells.Item(8, 9).Value = [double]"0.1271670650006795"
$ws.Cells.Item(8, 10).Value = [double]"0.1271670650006795"
$ws.Cells.Item(8, 13).Value = [double]"0.106124"
$ws.Cells.Item(8, 14).Value = [double]"0.318372"
$ws.Cells.Item(8, 15).Value = [double]"0.08094716512538251"
$ws.Cells.Item(8, 16).Value = [double]"0.08094716512538253"
$ws.Cells.Item(8, 17).Value = [double]"8.942372422193333"
$ws.Cells.Item(8, 18).Value = [double]"80.48135179974"
$ws.Cells.Item(8, 19).Value = [double]"0.01029381340912025"
$ws.Cells.Item(8, 20).Value = [double]"0.01029381340912025"
$ws.Cells.Item(9, 7).Value = [double]"84.26343166666668"
$ws.Cells.Item(9, 8).Value = [double]"252.790295"
$ws.Cells.Item(9, 9).Value = [double]"0.1271670650006795"
$ws.Cells.Item(9, 10).Value = [double]"0.1271670650006795"
$ws.Cells.Item(9, 15).Value = [double]"0.8331551016962769"
$ws.Cells.Item(9, 16).Value = [double]"0.833155101696277"
$ws.Cells.Item(9, 17).Value = [double]"92.04007568737278"
$ws.Cells.Item(9, 18).Value = [double]"828.360681186355"
$ws.Cells.Item(9, 19).Value = [double]"0.1059498889730582"
$ws.Cells.Item(9, 20).Value = [double]"0.1059498889730582"
$ws.Cells.Item(10, 7).Value = [double]"84.26343166666668"
$ws.Cells.Item(10, 8).Value = [double]"252.790295"
$ws.Cells.Item(10, 9).Value = [double]"0.1271670650006795"
$ws.Cells.Item(10, 10).Value = [double]"0.1271670650006795"
$ws.Cells.Item(10, 13).Value = [double]"0.1126143333333333"
$ws.Cells.Item(10, 14).Value = [double]"0.337843"
$ws.Cells.Item(10, 15).Value = [double]"0.08589773317834044"
$ws.Cells.Item(10, 16).Value = [double]"0.08589773317834046"
$ws.Cells.Item(10, 17).Value = [double]"9.489270181520556"
$ws.Cells.Item(10, 18).Value = [double]"85.403431633685"
$ws.Cells.Item(10, 19).Value = [double]"0.01092336261850104"
$ws.Cells.Item(10, 20).Value = [double]"0.01092336261850104"
$ws.Cells.Item(11, 5).Value = [double]"2"
$ws.Cells.Item(11, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(11, 7).Value = [double]"0.244783"
$ws.Cells.Item(11, 8).Value = [double]"0.7343489999999999"
$ws.Cells.Item(11, 9).Value = [double]"0.0003694168995537743"
$ws.Cells.Item(11, 10).Value = [double]"0.0003694168995537743"
$ws.Cells.Item(11, 13).Value = [double]"0.106124"
$ws.Cells.Item(11, 14).Value = [double]"0.318372"
$ws.Cells.Item(11, 15).Value = [double]"0.08094716512538251"
$ws.Cells.Item(11, 16).Value = [double]"0.08094716512538253"
$ws.Cells.Item(11, 17).Value = [double]"0.025977351092"
$ws.Cells.Item(11, 18).Value = [double]"0.233796159828"
$ws.Cells.Item(11, 19).Value = [double]"2.990325076828621E-05"
$ws.Cells.Item(11, 20).Value = [double]"2.990325076828621E-05"
$ws.Cells.Item(12, 5).Value = [double]"2"
$ws.Cells.Item(12, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(12, 7).Value = [double]"0.244783"
$ws.Cells.Item(12, 8).Value = [double]"0.7343489999999999"
$ws.Cells.Item(12, 9).Value = [double]"0.0003694168995537743"
$ws.Cells.Item(12, 10).Value = [double]"0.0003694168995537743"
$ws.Cells.Item(12, 15).Value = [double]"0.8331551016962769"
$ws.Cells.Item(12, 16).Value = [double]"0.833155101696277"
$ws.Cells.Item(12, 17).Value = [double]"0.2673739414756666"
$ws.Cells.Item(12, 18).Value = [double]"2.406365473281"
$ws.Cells.Item(12, 19).Value = [double]"0.0003077815745160481"
$ws.Cells.Item(12, 20).Value = [double]"0.0003077815745160482"
$ws.Cells.Item(13, 5).Value = [double]"2"
$ws.Cells.Item(13, 6).Value = [double]"0.6666666666666666"
$ws.Cells.Item(13, 7).Value = [double]"0.244783"
$ws.Cells.Item(13, 8).Value = [double]"0.7343489999999999"
$ws.Cells.Item(13, 9).Value = [double]"0.0003694168995537743"
$ws.Cells.Item(13, 10).Value = [double]"0.0003694168995537743"
$ws.Cells.Item(13, 13).Value = [double]"0.1126143333333333"
$ws.Cells.Item(13, 14).Value = [double]"0.337843"
$ws.Cells.Item(13, 15).Value = [double]"0.08589773317834044"
$ws.Cells.Item(13, 16).Value = [double]"0.08589773317834046"
$ws.Cells.Item(13, 17).Value = [double]"0.02756607435633333"
$ws.Cells.Item(13, 18).Value = [double]"0.248094669207"
$ws.Cells.Item(13, 19).Value = [double]"3.17320742694399E-05"
$ws.Cells.Item(13, 20).Value = [double]"3.17320742694399E-05"
$ws.Cells.Item(14, 7).Value = [double]"62.79827133333333"
$ws.Cells.Item(14, 8).Value = [double]"188.394814"
$ws.Cells.Item(14, 9).Value = [double]"0.09477268720988248"
$ws.Cells.Item(14, 10).Value = [double]"0.09477268720988248"
$ws.Cells.Item(14, 13).Value = [double]"0.106124"
$ws.Cells.Item(14, 14).Value = [double]"0.318372"
$ws.Cells.Item(14, 15).Value = [double]"0.08094716512538251"
$ws.Cells.Item(14, 16).Value = [double]"0.08094716512538253"
$ws.Cells.Item(14, 17).Value = [double]"6.664403746978667"
$ws.Cells.Item(14, 18).Value = [double]"59.979633722808"
$ws.Cells.Item(14, 19).Value = [double]"0.007671580360954585"
$ws.Cells.Item(14, 20).Value = [double]"0.007671580360954586"
$ws.Cells.Item(15, 7).Value = [double]"62.79827133333333"
$ws.Cells.Item(15, 8).Value = [double]"188.394814"
$ws.Cells.Item(15, 9).Value = [double]"0.09477268720988248"
$ws.Cells.Item(15, 10).Value = [double]"0.09477268720988248"
$ws.Cells.Item(15, 15).Value = [double]"0.8331551016962769"
$ws.Cells.Item(15, 16).Value = [double]"0.833155101696277"
$ws.Cells.Item(15, 17).Value = [double]"68.59390286192955"
$ws.Cells.Item(15, 18).Value = [double]"617.345125757366"
$ws.Cells.Item(15, 19).Value = [double]"0.07896034785037909"
$ws.Cells.Item(15, 20).Value = [double]"0.07896034785037909"
$ws.Cells.Item(16, 7).Value = [double]"62.79827133333333"
$ws.Cells.Item(16, 8).Value = [double]"188.394814"
$ws.Cells.Item(16, 9).Value = [double]"0.09477268720988248"
$ws.Cells.Item(16, 10).Value = [double]"0.09477268720988248"
$ws.Cells.Item(16, 13).Value = [double]"0.1126143333333333"
$ws.Cells.Item(16, 14).Value = [double]"0.337843"
$ws.Cells.Item(16, 15).Value = [double]"0.08589773317834044"
$ws.Cells.Item(16, 16).Value = [double]"0.08589773317834046"
$ws.Cells.Item(16, 17).Value = [double]"7.071985460689111"
$ws.Cells.Item(16, 18).Value = [double]"63.647869146202"
$ws.Cells.Item(16, 19).Value = [double]"0.008140758998548804"
$ws.Cells.Item(16, 20).Value = [double]"0.008140758998548804"
$ws.Cells.Item(17, 7).Value = [double]"512.5237530000001"
$ws.Cells.Item(17, 8).Value = [double]"1537.571259"
$ws.Cells.Item(17, 9).Value = [double]"0.7734807391890958"
$ws.Cells.Item(17, 10).Value = [double]"0.7734807391890958"
$ws.Cells.Item(17, 13).Value = [double]"0.106124"
$ws.Cells.Item(17, 14).Value = [double]"0.318372"
$ws.Cells.Item(17, 15).Value = [double]"0.08094716512538251"
$ws.Cells.Item(17, 16).Value = [double]"0.08094716512538253"
$ws.Cells.Item(17, 17).Value = [double]"54.39107076337201"
$ws.Cells.Item(17, 18).Value = [double]"489.519636870348"
$ws.Cells.Item(17, 19).Value = [double]"0.06261107311644266"
$ws.Cells.Item(17, 20).Value = [double]"0.06261107311644266"
$ws.Cells.Item(18, 7).Value = [double]"512.5237530000001"
$ws.Cells.Item(18, 8).Value = [double]"1537.571259"
$ws.Cells.Item(18, 9).Value = [double]"0.7734807391890958"
$ws.Cells.Item(18, 10).Value = [double]"0.7734807391890958"
$ws.Cells.Item(18, 15).Value = [double]"0.8331551016962769"
$ws.Cells.Item(18, 16).Value = [double]"0.833155101696277"
$ws.Cells.Item(18, 17).Value = [double]"559.824399323119"
$ws.Cells.Item(18, 18).Value = [double]"5038.419593908071"
$ws.Cells.Item(18, 19).Value = [double]"0.6444294239192025"
$ws.Cells.Item(18, 20).Value = [double]"0.6444294239192027"
$ws.Cells.Item(19, 7).Value = [double]"512.5237530000001"
$ws.Cells.Item(19, 8).Value = [double]"1537.571259"
$ws.Cells.Item(19, 9).Value = [double]"0.7734807391890958"
$ws.Cells.Item(19, 10).Value = [double]"0.7734807391890958"
$ws.Cells.Item(19, 13).Value = [double]"0.1126143333333333"
$ws.Cells.Item(19, 14).Value = [double]"0.337843"
$ws.Cells.Item(19, 15).Value = [double]"0.08589773317834044"
$ws.Cells.Item(19, 16).Value = [double]"0.08589773317834046"
$ws.Cells.Item(19, 17).Value = [double]"57.71752076159301"
$ws.Cells.Item(19, 18).Value = [double]"519.457686854337"
$ws.Cells.Item(19, 19).Value = [double]"0.06644024215345049"
$ws.Cells.Item(19, 20).Value = [double]"0.0664402421534505"

Write-Host "Applied TPM updates"
